# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the "全部类型" sheet to reflect the newly scraped data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 645
$ws1.Range("F4").Value = 1462
$ws1.Range("F5").Value = 684

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 645
$ws4.Range("F4").Value = 1462
$ws4.Range("F6").Value = 684
